# Results from June 11, 2020 12:01 AM run.
# The COVID disparities scraper re-ran: more states now have data (Alabama,
# Arkansas, California, Colorado), the existing rows were refreshed with
# the latest published numbers, and the whole table is now sorted
# alphabetically by Location (it previously was not sorted). Rebuild the
# data block (rows 2-16) from scratch to land on the new row order/values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data region (rows 2-12) completely; we will rewrite
# rows 2-16 from scratch in the new (alphabetically-sorted) order.
$ws.Range("A2:I16").Clear()

# Row 2: Alabama
$ws.Range("A2").Value = "Alabama"
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"
$ws.Range("B2").Value = 43992
$ws.Range("C2").Value = 21626
$ws.Range("D2").Value = 739
$ws.Range("E2").Value = 9221
$ws.Range("F2").Value = 333
$ws.Range("G2").Value = 42.64
$ws.Range("H2").Value = 45.06
$ws.Range("I2").Value = "Success!"

# Row 3: Arkansas
$ws.Range("A3").Value = "Arkansas"
$ws.Range("B3").NumberFormat = "YYYY-MM-DD"
$ws.Range("B3").Value = 43992
$ws.Range("C3").Value = 593
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 0
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2.7"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "0.0"
$ws.Range("I3").Value = "Success!"

# Row 4: California
$ws.Range("A4").Value = "California"
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"
$ws.Range("B4").Value = 43991
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "97336"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4600"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4713"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "451"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "4.8"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "9.8"
$ws.Range("I4").Value = "Success!"

# Row 5: California - San Diego
$ws.Range("A5").Value = "California - San Diego"
$ws.Range("B5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B5").Value = 43991.99836937636
$ws.Range("C5").Value = 8729
$ws.Range("D5").Value = 301
$ws.Range("E5").Value = 260
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 2.98
$ws.Range("H5").Value = 2.99
$ws.Range("I5").Value = "Success!"

# Row 6: Colorado
$ws.Range("A6").Value = "Colorado"
$ws.Range("B6").NumberFormat = "YYYY-MM-DD"
$ws.Range("B6").Value = 43992
$ws.Range("C6").Value = 28499
$ws.Range("D6").Value = 1573
$ws.Range("E6").Value = 1676
$ws.Range("F6").Value = 106
$ws.Range("G6").Value = 5.88
$ws.Range("H6").Value = 6.74
$ws.Range("I6").Value = "Success!"

# Row 7: Florida
$ws.Range("A7").Value = "Florida"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2020-06-10"
$ws.Range("C7").Value = 65779
$ws.Range("D7").Value = 2801
$ws.Range("E7").Value = 12198
$ws.Range("F7").Value = 558
$ws.Range("G7").Value = 18.54
$ws.Range("H7").Value = 19.92
$ws.Range("I7").Value = "Success!"

# Row 8: Georgia
$ws.Range("A8").Value = "Georgia"
$ws.Range("B8").NumberFormat = "YYYY-MM-DD"
$ws.Range("B8").Value = 43992
$ws.Range("C8").Value = 53980
$ws.Range("D8").Value = 2329
$ws.Range("E8").Value = 16965
$ws.Range("F8").Value = 1123
$ws.Range("G8").Value = 31.43
$ws.Range("H8").Value = 48.22
$ws.Range("I8").Value = "Success!"

# Row 9: Massachusetts
$ws.Range("A9").Value = "Massachusetts"
$ws.Range("B9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B9").Value = 43992
$ws.Range("C9").Value = 104156
$ws.Range("D9").Value = 7454
$ws.Range("E9").Value = 9729
$ws.Range("F9").Value = 620
$ws.Range("G9").Value = 9.34
$ws.Range("H9").Value = 8.32
$ws.Range("I9").Value = "Success!"

# Row 10: Michigan
$ws.Range("A10").Value = "Michigan"
$ws.Range("B10:H10").Value = ""
$ws.Range("I10").Value = "An error occurred. ... UnboundLocalError(`"local variable 'date_published' referenced before assignment`")"

# Row 11: Minnesota
$ws.Range("A11").Value = "Minnesota"
$ws.Range("B11").NumberFormat = "YYYY-MM-DD"
$ws.Range("B11").Value = 43992
$ws.Range("C11").Value = 28869
$ws.Range("D11").Value = 1236
$ws.Range("E11").Value = 6342
$ws.Range("F11").Value = 78
$ws.Range("G11").Value = 21.97
$ws.Range("H11").Value = 6.31
$ws.Range("I11").Value = "Success!"

# Row 12: North Carolina
$ws.Range("A12").Value = "North Carolina"
$ws.Range("B12:H12").Value = ""
$ws.Range("I12").Value = "An error occurred. ... ValueError('Unable to extract date from table header.')"

# Row 13: Texas -- Bexar County
$ws.Range("A13").Value = "Texas -- Bexar County"
$ws.Range("B13").NumberFormat = "YYYY-MM-DD"
$ws.Range("B13").Value = 43992
$ws.Range("C13").Value = 1805
$ws.Range("D13").Value = 54
$ws.Range("E13").Value = 280
$ws.Range("F13").Value = 17
$ws.Range("G13").Value = 15.51
$ws.Range("H13").Value = 31.48
$ws.Range("I13").Value = "Success!"

# Row 14: Virginia
$ws.Range("A14").Value = "Virginia"
$ws.Range("B14:H14").Value = ""
$ws.Range("I14").Value = "An error occurred. ... URLError(TimeoutError(10060, 'A connection attempt failed because the connected party did not properly respond after a period of time, or established connection failed because connected host has failed to respond', None, 10060, None))"

# Row 15: Washington, DC
$ws.Range("A15").Value = "Washington, DC"
$ws.Range("B15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B15").Value = 43991
$ws.Range("C15").Value = 9474
$ws.Range("D15").Value = 495
$ws.Range("E15").Value = 4331
$ws.Range("F15").Value = 367
$ws.Range("G15").Value = 45.71
$ws.Range("H15").Value = 74.14
$ws.Range("I15").Value = "Success!"

# Row 16: Wisconsin -- Milwaukee
$ws.Range("A16").Value = "Wisconsin -- Milwaukee"
$ws.Range("B16").NumberFormat = "YYYY-MM-DD"
$ws.Range("B16").Value = 43992
$ws.Range("C16").Value = 9161
$ws.Range("D16").Value = 306
$ws.Range("E16").Value = 2597
$ws.Range("F16").Value = 129
$ws.Range("G16").Value = 28.35
$ws.Range("H16").Value = 1.41
$ws.Range("I16").Value = "Success!"

